$wb = $excel.ActiveWorkbook

# Sheet order in the workbook:
# 1 Funciones_Objetivo
# 2 Restricciones_del_lider
# 3 Restricciones_del_follower
# 4 Punto_modificado
# 5 Vector_bf
# 6 Vector_BF
# 7 Vector_Alpha
# Use numeric indices to avoid case-insensitive name collisions (Vector_bf vs Vector_BF).

# Helper: the source file stores these numeric-looking values as TEXT
# (shared strings), not as numbers. Excel's normal Value assignment
# auto-converts a numeric-looking string into a real number, so we
# momentarily force a Text number format, assign the value, then clear
# the format again so the cell keeps its original (default) appearance
# while the stored cell type remains text.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$wsFollower = $wb.Worksheets.Item(3)
Set-TextValue $wsFollower.Range("A2") "-34.13048780487805 - 2x_1 + 7.219512195121951y_1 + 3.975609756097562y_2"
Set-TextValue $wsFollower.Range("B2") "36.63048780487805"
Set-TextValue $wsFollower.Range("D2") "0.96"
Set-TextValue $wsFollower.Range("E2") "0.8"
Set-TextValue $wsFollower.Range("F2") "7.800000000000001"

Set-TextValue $wsFollower.Range("A3") "1.8470720754249346 + x_1 - 3x_2 - 0.4685212856658182y_1 + 0.6251829714673454y_2"
Set-TextValue $wsFollower.Range("B3") "-3.8470720754249346"
Set-TextValue $wsFollower.Range("D3") "0.9"
Set-TextValue $wsFollower.Range("E3") "8.100000000000001"
Set-TextValue $wsFollower.Range("F3") "0"

Set-TextValue $wsFollower.Range("A4") "-3.9499999999999997 + x_1 + x_2"
Set-TextValue $wsFollower.Range("B4") "1.7999999999999998"
Set-TextValue $wsFollower.Range("D4") "0.28"
Set-TextValue $wsFollower.Range("E4") "6.5"
Set-TextValue $wsFollower.Range("F4") "0"

$wsPunto = $wb.Worksheets.Item(4)
Set-TextValue $wsPunto.Range("A2") "2.55"
Set-TextValue $wsPunto.Range("B2") "1.25"
Set-TextValue $wsPunto.Range("C2") "4.25"
Set-TextValue $wsPunto.Range("D2") "2.15"

$wsbf = $wb.Worksheets.Item(5)
Set-TextValue $wsbf.Range("A2") "-2.5090625502178368"
Set-TextValue $wsbf.Range("A3") "-5.37925004017427"

$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF.Range("A2") "-11.000000000000002"
Set-TextValue $wsBF.Range("A3") "16.800000000000004"
Set-TextValue $wsBF.Range("A4") "-2.4805873422044336"
Set-TextValue $wsBF.Range("A5") "-8.24446987376355"

# Vector_Alpha stores its numeric values as real numbers (not text), so
# assign them directly as numbers.
$wsAlpha = $wb.Worksheets.Item(7)
$wsAlpha.Range("A2").Value = 0.75
$wsAlpha.Range("A3").Value = 0.6000000000000001
